# UndoRedoActivityDiagram.pptx edit
#
# 1. Refresh the cached text of the auto-updating "datetimeFigureOut"
#    date field that appears on the slide master and every slide layout
#    (it was re-saved on a later date, so the cached text moved from
#    30/12/2017 to 16/4/2018).
# 2. On the diagram slide: widen the "Add command to undo stack" box and
#    retarget its wording to "Add addressbook to undo stack" (the word
#    "command" was replaced by "addressbook").
# 3. Re-anchor the elbow connector that runs from that box so it keeps
#    meeting the box's (now-moved) right-hand connection point.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name.StartsWith("Date Placeholder")) {
            $shp.TextFrame.TextRange.Text = "16/4/2018"
        }
    }
}

# --- 1. Slide master + every slide layout: refresh the date placeholder ---
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholder $layouts.Item($L).Shapes
}

# --- 2. Diagram slide: widen the rounded rectangle & update its text ---
$s = $p.Slides.Item(1)

$rect = $null
$connector = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Rounded Rectangle 50") {
        $rect = $shp
    }
    if ($shp.Name -eq "Elbow Connector 73") {
        $connector = $shp
    }
}

# Widen the box: cx 1634410 EMU -> 1870298 EMU (914400 EMU = 1 inch = 72 pt).
# Shape geometry is exposed in points (single-precision) over COM, so a
# tiny epsilon compensates for the EMU<->pt round trip truncating down.
$emuPerPt = 12700
$epsilon = 0.00005
$rect.Width = (1870298 / $emuPerPt) + $epsilon

# Replace the word "command" with "addressbook" in place (keeps the
# surrounding text, mirroring how PowerPoint splits the run when a user
# retypes a single word in the middle of a sentence).
$tr = $rect.TextFrame.TextRange
$commandStart = $tr.Text.IndexOf("command") + 1
$sub = $tr.Characters($commandStart, "command".Length)
$sub.Text = "addressbook"

# Also nudge the trailing " to undo stack" so it mirrors the same
# resulting run boundary as the source edit ("  to  " | "undo stack").
$afterText = $tr.Text
$undoStart = $afterText.IndexOf("undo stack") + 1
$sub2 = $tr.Characters($undoStart, "undo stack".Length)
$sub2.Text = "undo stack"

# --- 3. Re-anchor the connector leaving the box's right edge ---
$connector.Left = (9058496 / $emuPerPt) + $epsilon
$connector.Width = (240382 / $emuPerPt) + $epsilon
